$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (duplicate_image_filename) should be filled with "NA" for the
# data rows 2 through 21 (previously blank in that column).
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
